# Apply "added notes and other images" edit:
#  - recolor a few topic rows (A16, A21, A22 -> green; A23 -> red)
#  - add two "notes" cells in column B (B24, B25) referencing the Guru 99
#    resource, with the same orange fill used by the existing notes column
#  - move the view/selection down to the newly edited area (A24)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Recolor existing topic cells -----------------------------------------
# Green fill (matches existing "00B050" fill already used elsewhere)
$ws.Range("A16").Interior.Color = 5287936
$ws.Range("A21").Interior.Color = 5287936
$ws.Range("A22").Interior.Color = 5287936

# Red fill (matches existing "FF0000" fill already used elsewhere)
$ws.Range("A23").Interior.Color = 255

# --- Add new notes in column B ---------------------------------------------
# Write B25 before B24 so the new shared-string entries land in the same
# order as the source edit (index 32 = long note, index 33 = short note).
$ws.Range("B25").Value = "Guru 99 - WBS, FP method, 3 point estimation"
$ws.Range("B25").Interior.Color = 49407

$ws.Range("B24").Value = "Guru 99"
$ws.Range("B24").Interior.Color = 49407

# --- Update the view/selection ---------------------------------------------
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 21
$aw.ScrollColumn = 1
$ws.Range("A24").Select()
